# Auto-generated edit script applying numeric corrections to Leve profit sheets
# per commit 'chore: update Sheets via scheduled runner'
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Cells.Item(18, 8).Value = 2100
$ws.Cells.Item(18, 9).Value = 2100
$ws.Cells.Item(18, 11).Value = 2100
$ws.Cells.Item(18, 13).Value = -1816

# Row 51
$ws.Cells.Item(51, 8).Value = 5000
$ws.Cells.Item(51, 9).Value = 5000
$ws.Cells.Item(51, 11).Value = 5000
$ws.Cells.Item(51, 13).Value = -4516

# Row 88
$ws.Cells.Item(88, 8).Value = 1763.4546
$ws.Cells.Item(88, 9).Value = 1056.25
$ws.Cells.Item(88, 10).Value = 2167.5715
$ws.Cells.Item(88, 11).Value = 1056.25
$ws.Cells.Item(88, 12).Value = 2167.5715
$ws.Cells.Item(88, 13).Value = -650.25
$ws.Cells.Item(88, 14).Value = -2979.5715

# Row 91
$ws.Cells.Item(91, 8).Value = 1763.4546
$ws.Cells.Item(91, 9).Value = 1056.25
$ws.Cells.Item(91, 10).Value = 2167.5715
$ws.Cells.Item(91, 11).Value = 1056.25
$ws.Cells.Item(91, 12).Value = 2167.5715
$ws.Cells.Item(91, 13).Value = 347.75
$ws.Cells.Item(91, 14).Value = -4975.5715

# Row 96
$ws.Cells.Item(96, 8).Value = 1598.125
$ws.Cells.Item(96, 9).Value = 1076.8182
$ws.Cells.Item(96, 10).Value = 2745
$ws.Cells.Item(96, 11).Value = 3230.4546
$ws.Cells.Item(96, 12).Value = 8235
$ws.Cells.Item(96, 13).Value = -1857.4546
$ws.Cells.Item(96, 14).Value = -10981

# Row 137
$ws.Cells.Item(137, 8).Value = 4604.1665
$ws.Cells.Item(137, 10).Value = 7700.4546
$ws.Cells.Item(137, 12).Value = 23101.3638
$ws.Cells.Item(137, 14).Value = -28201.3638

# Row 138
$ws.Cells.Item(138, 8).Value = 6162.212
$ws.Cells.Item(138, 9).Value = 6377.316
$ws.Cells.Item(138, 10).Value = 5870.2856
$ws.Cells.Item(138, 11).Value = 19131.948
$ws.Cells.Item(138, 12).Value = 17610.8568
$ws.Cells.Item(138, 13).Value = -13991.948
$ws.Cells.Item(138, 14).Value = -27890.8568

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Cells.Item(61, 8).Value = 2239.6
$ws.Cells.Item(61, 9).Value = 2239.6
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 2239.6
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -2027.6
$ws.Cells.Item(61, 14).ClearContents()

# Row 122
$ws.Cells.Item(122, 8).Value = 3877.2222
$ws.Cells.Item(122, 9).Value = 3104.875
$ws.Cells.Item(122, 10).Value = 4495.1
$ws.Cells.Item(122, 11).Value = 9314.625
$ws.Cells.Item(122, 12).Value = 13485.3
$ws.Cells.Item(122, 13).Value = -6864.625
$ws.Cells.Item(122, 14).Value = -18385.3

# Row 136
$ws.Cells.Item(136, 8).Value = 2239.6
$ws.Cells.Item(136, 9).Value = 2239.6
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 6718.799999999999
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -4168.799999999999
$ws.Cells.Item(136, 14).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 96
$ws.Cells.Item(96, 8).Value = 11250
$ws.Cells.Item(96, 9).Value = 11250
$ws.Cells.Item(96, 11).Value = 11250
$ws.Cells.Item(96, 13).Value = -8504

# Row 134
$ws.Cells.Item(134, 8).Value = 2486.9285
$ws.Cells.Item(134, 10).Value = 4798.1665
$ws.Cells.Item(134, 12).Value = 14394.4995
$ws.Cells.Item(134, 14).Value = -19464.4995

$ws = $wb.Worksheets.Item("CRP")
# Row 21
$ws.Cells.Item(21, 8).Value = 7777
$ws.Cells.Item(21, 10).Value = 7777
$ws.Cells.Item(21, 12).Value = 7777
$ws.Cells.Item(21, 14).Value = -8247

# Row 22
$ws.Cells.Item(22, 8).Value = 324.5
$ws.Cells.Item(22, 9).Value = 149
$ws.Cells.Item(22, 11).Value = 149
$ws.Cells.Item(22, 13).Value = 201

# Row 58
$ws.Cells.Item(58, 8).Value = 5558
$ws.Cells.Item(58, 9).Value = 1891.8
$ws.Cells.Item(58, 11).Value = 1891.8
$ws.Cells.Item(58, 13).Value = -1688.8

# Row 99
$ws.Cells.Item(99, 8).Value = 16736.455
$ws.Cells.Item(99, 9).Value = 15094.308
$ws.Cells.Item(99, 10).Value = 19108.445
$ws.Cells.Item(99, 11).Value = 15094.308
$ws.Cells.Item(99, 12).Value = 19108.445
$ws.Cells.Item(99, 13).Value = -13596.308
$ws.Cells.Item(99, 14).Value = -22104.445

# Row 122
$ws.Cells.Item(122, 8).Value = 909
$ws.Cells.Item(122, 9).Value = 886.25
$ws.Cells.Item(122, 11).Value = 2658.75
$ws.Cells.Item(122, 13).Value = -208.75

# Row 126
$ws.Cells.Item(126, 8).Value = 16736.455
$ws.Cells.Item(126, 9).Value = 15094.308
$ws.Cells.Item(126, 10).Value = 19108.445
$ws.Cells.Item(126, 11).Value = 45282.924
$ws.Cells.Item(126, 12).Value = 57325.335
$ws.Cells.Item(126, 13).Value = -42812.924
$ws.Cells.Item(126, 14).Value = -62265.335

# Row 132
$ws.Cells.Item(132, 8).Value = 1733.4584
$ws.Cells.Item(132, 9).Value = 1733.4584
$ws.Cells.Item(132, 11).Value = 5200.3752
$ws.Cells.Item(132, 13).Value = -2670.3752

# Row 136
$ws.Cells.Item(136, 8).Value = 5558
$ws.Cells.Item(136, 9).Value = 1891.8
$ws.Cells.Item(136, 11).Value = 5675.4
$ws.Cells.Item(136, 13).Value = -3125.4

$ws = $wb.Worksheets.Item("CUL")
# Row 21
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 14).ClearContents()

# Row 26
$ws.Cells.Item(26, 8).Value = 181.66667
$ws.Cells.Item(26, 9).Value = 181.66667
$ws.Cells.Item(26, 11).Value = 545.00001
$ws.Cells.Item(26, 13).Value = -257.00001

# Row 32
$ws.Cells.Item(32, 8).Value = 2911025.8
$ws.Cells.Item(32, 10).Value = 4365039
$ws.Cells.Item(32, 12).Value = 13095117
$ws.Cells.Item(32, 14).Value = -13095683

# Row 34
$ws.Cells.Item(34, 8).Value = 2017.1875
$ws.Cells.Item(34, 10).Value = 2512.5
$ws.Cells.Item(34, 12).Value = 7537.5
$ws.Cells.Item(34, 14).Value = -7705.5

# Row 63
$ws.Cells.Item(63, 8).Value = 1599
$ws.Cells.Item(63, 9).Value = 1599
$ws.Cells.Item(63, 11).Value = 4797
$ws.Cells.Item(63, 13).Value = -4048

# Row 66
$ws.Cells.Item(66, 8).Value = 1599
$ws.Cells.Item(66, 9).Value = 1599
$ws.Cells.Item(66, 11).Value = 14391
$ws.Cells.Item(66, 13).Value = -10647

# Row 107
$ws.Cells.Item(107, 8).Value = 3623.5
$ws.Cells.Item(107, 10).Value = 1748
$ws.Cells.Item(107, 12).Value = 5244
$ws.Cells.Item(107, 14).Value = -9084

# Row 129
$ws.Cells.Item(129, 8).Value = 2637.375
$ws.Cells.Item(129, 10).Value = 2454.75
$ws.Cells.Item(129, 12).Value = 7364.25
$ws.Cells.Item(129, 14).Value = -17364.25

# Row 137
$ws.Cells.Item(137, 8).Value = 3565.0715
$ws.Cells.Item(137, 9).Value = 3912.5
$ws.Cells.Item(137, 10).Value = 3426.1
$ws.Cells.Item(137, 11).Value = 11737.5
$ws.Cells.Item(137, 12).Value = 10278.3
$ws.Cells.Item(137, 13).Value = -6637.5
$ws.Cells.Item(137, 14).Value = -20478.3

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Cells.Item(122, 8).Value = 921827.7
$ws.Cells.Item(122, 9).Value = 204588.6
$ws.Cells.Item(122, 10).Value = 1434141.2
$ws.Cells.Item(122, 11).Value = 613765.8
$ws.Cells.Item(122, 12).Value = 4302423.6
$ws.Cells.Item(122, 13).Value = -611315.8
$ws.Cells.Item(122, 14).Value = -4307323.6

# Row 132
$ws.Cells.Item(132, 8).Value = 3769.5557
$ws.Cells.Item(132, 9).Value = 2432
$ws.Cells.Item(132, 11).Value = 7296
$ws.Cells.Item(132, 13).Value = -4766

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Cells.Item(55, 8).Value = 952.2222
$ws.Cells.Item(55, 9).Value = 674.8
$ws.Cells.Item(55, 10).Value = 1299
$ws.Cells.Item(55, 11).Value = 674.8
$ws.Cells.Item(55, 12).Value = 1299
$ws.Cells.Item(55, 13).Value = -501.8
$ws.Cells.Item(55, 14).Value = -1645

# Row 122
$ws.Cells.Item(122, 8).Value = 4166.3335
$ws.Cells.Item(122, 9).Value = 3999.6
$ws.Cells.Item(122, 10).Value = 5000
$ws.Cells.Item(122, 11).Value = 11998.8
$ws.Cells.Item(122, 12).Value = 15000
$ws.Cells.Item(122, 13).Value = -9548.799999999999
$ws.Cells.Item(122, 14).Value = -19900

# Row 132
$ws.Cells.Item(132, 8).Value = 6499.857
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 13).ClearContents()

# Row 136
$ws.Cells.Item(136, 8).Value = 4999
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 13).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Cells.Item(81, 8).Value = 2998.9092
$ws.Cells.Item(81, 9).Value = 2998.9092
$ws.Cells.Item(81, 11).Value = 5997.8184
$ws.Cells.Item(81, 13).Value = -4936.8184

# Row 84
$ws.Cells.Item(84, 8).Value = 2998.9092
$ws.Cells.Item(84, 9).Value = 2998.9092
$ws.Cells.Item(84, 11).Value = 29989.092
$ws.Cells.Item(84, 13).Value = -24685.092

# Row 122
$ws.Cells.Item(122, 8).Value = 10899.2
$ws.Cells.Item(122, 9).Value = 3124.25
$ws.Cells.Item(122, 11).Value = 9372.75
$ws.Cells.Item(122, 13).Value = -6922.75

# Row 126
$ws.Cells.Item(126, 8).Value = 65628.69
$ws.Cells.Item(126, 9).Value = 253123.75
$ws.Cells.Item(126, 10).Value = 3130.3333
$ws.Cells.Item(126, 11).Value = 759371.25
$ws.Cells.Item(126, 12).Value = 9390.999899999999
$ws.Cells.Item(126, 13).Value = -756901.25
$ws.Cells.Item(126, 14).Value = -14330.9999

# Row 132
$ws.Cells.Item(132, 8).Value = 2735.2104
$ws.Cells.Item(132, 9).Value = 2123
$ws.Cells.Item(132, 10).Value = 3784.7144
$ws.Cells.Item(132, 11).Value = 6369
$ws.Cells.Item(132, 12).Value = 11354.1432
$ws.Cells.Item(132, 13).Value = -3839
$ws.Cells.Item(132, 14).Value = -16414.1432
